$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Restyle existing row 65 (border style) and add new A65 value ---
$ws.Range("A4:E4").Copy()
$ws.Range("A65:E65").PasteSpecial(-4122)
$ws.Range("A65").Value = 'SCRIPT/T01P02A/um2504.ssb'
$ws.Rows.Item(65).RowHeight = 49.2

# --- Row 66: normal style, A present ---
$ws.Range("A61:E61").Copy()
$ws.Range("A66:E66").PasteSpecial(-4122)
$ws.Range("A66").Value = 'SCRIPT/T01P02A/us0101.ssb'
$ws.Range("B66").Value = 358
$ws.Range("C66").Value = ' Yo! I heard what\''s going on with\nyou two!'
$ws.Range("D66").Value = ' Йо! Я тут о вас двоих кое-что\nузнал!'
$ws.Range("E66").Value = ' Êï! Ÿ óôó ï âàò äâïéö ëïå-œóï\nôèîàì!'
$ws.Rows.Item(66).RowHeight = 43.2

# --- Row 67: border style, A present-empty ---
$ws.Range("A60:E60").Copy()
$ws.Range("A67:E67").PasteSpecial(-4122)
$ws.Range("B67").Value = 361
$ws.Range("C67").Value = ' You\''re taking the guild\''s\ngraduation exam? You give it your all!'
$ws.Range("D67").Value = ' Вы собираетесь пройти выпускной\nэкзамен? Выложитесь на полную!'
$ws.Range("E67").Value = ' Âú òïáéñàåóåòû ðñïêóé âúðôòëîïê\nüëèàíåî? Âúìïçéóåòû îà ðïìîôý!'
$ws.Rows.Item(67).RowHeight = 31.8

# --- Row 68: normal style, A present ---
$ws.Range("A61:E61").Copy()
$ws.Range("A68:E68").PasteSpecial(-4122)
$ws.Range("A68").Value = 'SCRIPT/T01P02A/us3104.ssb'
$ws.Range("B68").Value = 333
$ws.Range("C68").Value = ' Yo! You passed your graduation\nexam, did you?'
$ws.Range("D68").Value = ' Йо! Вы прошли выпускной экзамен,\nда?'
$ws.Range("E68").Value = ' Êï! Âú ðñïšìé âúðôòëîïê üëèàíåî,\näà?'
$ws.Rows.Item(68).RowHeight = 43.2

# --- Row 69: normal style, NO A cell ---
$ws.Range("B62:E62").Copy()
$ws.Range("B69:E69").PasteSpecial(-4122)
$ws.Range("B69").Value = 336
$ws.Range("C69").Value = ' You did it! Congratulations!'
$ws.Range("D69").Value = ' Молодцы! Поздравляю!'
$ws.Range("E69").Value = ' Íïìïäøú! Ðïèäñàâìÿý!'

# --- Row 70: normal style, NO A cell ---
$ws.Range("B62:E62").Copy()
$ws.Range("B70:E70").PasteSpecial(-4122)
$ws.Range("B70").Value = 339
$ws.Range("C70").Value = ' I recognized that you two\nwere something else. And I was right!'
$ws.Range("D70").Value = ' Я сразу понял, что вы справитесь.\nИ я не ошибся!'
$ws.Range("E70").Value = ' Ÿ òñàèô ðïîÿì, œóï âú òðñàâéóåòû.\nÉ ÿ îå ïšéáòÿ!'
$ws.Rows.Item(70).RowHeight = 21.6

# --- Update view: selection D72 ---
$ws.Range("D72").Select()

